# Translate the "Gesamtübersicht" (overview) report-template sheet from
# German example text to English, matching the commit's "translated the
# example" change. New strings are appended to the shared-string table in
# the same order as the cells are written below (A2, A3, B3, D3..J3, A7,
# A8, A9) so they land at shared-string indices 27-39, same as target.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: section heading
$ws.Range("A2").Value2 = "Accounting period"

# Row 3: column headers
$ws.Range("A3").Value2 = "From"
$ws.Range("B3").Value2 = "Until"
# C3 "Budget" is unchanged
$ws.Range("D3").Value2 = "Spent net"
$ws.Range("E3").Value2 = "Spent gross"
$ws.Range("F3").Value2 = "Hours"
$ws.Range("G3").Value2 = "Remaining budget net"
$ws.Range("H3").Value2 = "Remaining budget gross"
$ws.Range("I3").Value2 = "Progress"
$ws.Range("J3").Value2 = "Invoice recipient"

# Row 7: sub-total label
$ws.Range("A7").Value2 = "Sum according to invoice recipient"

# Row 8: per-name sum label
$ws.Range("A8").Value2 = "Sum {name}"

# Row 9: grand total label
$ws.Range("A9").Value2 = "Total sum"

# The sheet's cursor/selection ends up on A14 after the edits.
$ws.Range("A14").Select() | Out-Null
